# "reg file fixed bugs and first CB implementation"
# Adds a new "status" tracking column to the work plan and switches the
# sheet's font from Calibri to Arial.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D: status tracking, with the first task marked "done"
$ws.Range("D1").Value = "status"
$ws.Range("D2").Value = "done"

# Switch the whole (now A1:D16) sheet over to Arial
$ws.Range("A1:D16").Font.Name = "Arial"

# Leave the selection where the author ended up after typing "done"
$ws.Range("D3").Select() | Out-Null
